# Update the "Förändrad" (Changed) date column C for rows 2-11
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C11").Value = 45175
